$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================================
# 1) Re-purpose the style that currently lives on C2:C7 (wrap, vertical-top)
#    for the brand-new E:G columns *before* we touch C's own styling, since
#    once C adopts B's style the s=4 "donor" combo would have to be rebuilt.
# =========================================================================
$ws.Range("C2").Copy()
$ws.Range("E2:G7").PasteSpecial(-4122)

# =========================================================================
# 2) Cell values
# =========================================================================

# ---- Header row ----
$ws.Range("A1").Value = 'codice_ 1_livello'
$ws.Range("B1").Value = 'label_ITA_1_Livello'
$ws.Range("C1").Value = 'label_ENG_1_livello'
$ws.Range("D1").Value = 'label_DEU_1_livello'
$ws.Range("E1").Value = 'definizione_ITA'
$ws.Range("F1").Value = 'definizione_ENG'
$ws.Range("G1").Value = 'definizione_DEU'

# ---- Row 2 ----
$ws.Range("A2").Value = 'NONE'
$ws.Range("B2").Value = 'Non online'
$ws.Range("C2").Value = 'Not Online'
$ws.Range("D2").Value = 'Nicht online'
$ws.Range("E2").Value = 'Servizio erogato offline, per il quale non esiste nemmeno una pagina web informativa.'
$ws.Range("F2").Value = 'Offline service for which it does not exist any informational web page'
$ws.Range("G2").Value = 'Offline bereitgestellter Dienst, für den es nicht einmal eine informative Webseite gibt.'

# ---- Row 3 ----
$ws.Range("A3").Value = 'LEVEL 1'
$ws.Range("B3").Value = 'Informazione'
$ws.Range("C3").Value = 'Information'
$ws.Range("D3").Value = 'Informationen'
$ws.Range("E3").Value = 'Sono fornite all''utente informazioni sul procedimento amministrativo (es. finalità, termini di richiesta, ecc.) ed eventualmente sulle modalità di espletamento (es. sedi, orari di sportello).'
$ws.Range("F3").Value = 'Users are informed about the administrative process that regards the service (e.g., objectives, how to require it, etc) and about the way in which the service can be used (e.g., opening hours of the information desk, location, etc)'
$ws.Range("G3").Value = 'Der Benutzer erhält Informationen über den administrativen Ablauf (z. B. Zweck, Bedingungen der Anfrage usw.) und eventuell über die Art und Weise, wie er zu erledigen ist (z. B. Stellen, Schalterzeiten).'

# ---- Row 4 ----
$ws.Range("A4").Value = 'LEVEL 2'
$ws.Range("B4").Value = 'Interazione ad una via '
$ws.Range("C4").Value = 'One way interaction'
$ws.Range("D4").Value = 'Einweg-Interaktion'
$ws.Range("E4").Value = 'Oltre alle informazioni, sono resi disponibili all''utente i moduli per la richiesta dell''atto/procedimento amministrativo di interesse che dovrà poi essere inoltrata attraverso canali tradizionali (es. modulo di variazione residenza o moduli di autocertificazione).'
$ws.Range("F4").Value = 'In addition to the information, users can use online forms in order to start a request for an administrative act of interest. Afterwords, the user must send the forms via traditional channels'
$ws.Range("G4").Value = 'Zusätzlich zu den Informationen werden dem Benutzer die Formulare zur Anforderung des gewünschten Verwaltungsaktes / des Verfahrens zur Verfügung gestellt, die dann über herkömmliche Kanäle weitergeleitet werden muss (z. B. Formular zur Änderung der Residenz oder Formulare zur Selbstbescheinigung).'

# ---- Row 5 ----
$ws.Range("A5").Value = 'LEVEL 3'
$ws.Range("B5").Value = 'Interazione a due vie'
$ws.Range("C5").Value = 'Bidirectional interaction'
$ws.Range("D5").Value = 'Bidirektionale Interaktion'
$ws.Range("E5").Value = 'L''utente può avviare l''atto/procedimento amministrativo di interesse (es. il modulo può essere compilato e inviato on line) e viene garantita on line solo la presa in carico dei dati immessi dall''utente e non la loro contestuale elaborazione.'
$ws.Range("F5").Value = 'The user can start an administrative act of interest online (e.g., the form can be filled in and sent online) and it is guaranteed that the data is provided online, only; it is not guaranteed the concurrent online data processing.'
$ws.Range("G5").Value = 'Der Benutzer kann den gewünschten Verwaltungsakt / Vorgang initiieren (z. B. das Formular kann online ausgefüllt und gesendet werden) und es werden online nur die Übernahme der vom Benutzer eingegebenen Daten garantiert, nicht aber deren kontextuelle Verarbeitung.'

# ---- Row 6 ----
$ws.Range("A6").Value = 'LEVEL 4'
$ws.Range("B6").Value = 'Transazione'
$ws.Range("C6").Value = 'Transaction'
$ws.Range("D6").Value = 'Transaktion'
$ws.Range("E6").Value = 'L''utente può avviare l''atto/procedimento amministrativo di interesse fornendo i dati necessari ed eseguire la transazione corrispondente interamente on line, incluso l''eventuale pagamento dei costi previsti.'
$ws.Range("F6").Value = 'The user can start an administrative act of interest online by providing the necessary data. The user can then carry out the transaction entirely online, including the possible payment of costs related to the service.'
$ws.Range("G6").Value = 'Der Benutzer kann den gewünschten Verwaltungsakt / Vorgang initiieren, indem er die erforderlichen Daten zur Verfügung stellt und die entsprechende Transaktion vollständig online durchführt, einschließlich der Zahlung der erwarteten Kosten.'

# ---- Row 7 ----
$ws.Range("A7").Value = 'LEVEL 5'
$ws.Range("B7").Value = 'Personalizzazione'
$ws.Range("C7").Value = 'Customization'
$ws.Range("D7").Value = 'Personalisierung'
$ws.Range("E7").Value = 'L''utente, oltre ad eseguire on line l''intero ciclo del procedimento amministrativo di interesse riceve informazioni (sono ricordate le scadenze, è restituito l''esito del procedimento, ecc.), che gli sono inviate preventivamente, sulla base del profilo collegato (c.d. pro-attività).'
$ws.Range("F7").Value = 'The user can carried out an entire administrative act online and (s)he can receive information according to his/her profile (e.g., it can receive information about deadlines, about the result of an administrative act, etc.) This is also called pro-active service.'
$ws.Range("G7").Value = 'Zusätzlich zum gesamten Zyklus des Verwaltungsverfahrens von Interesse, den er selbst online ausführt, erhält der Benutzer Informationen (Erinnerung an die Fristen, Zurückgabe des Verfahrensergebnisses, etc.), die ihm auf der Besis des verknüpften Profils im Voraus gesendet werden (Pro-Aktivität)'

# =========================================================================
# 3) Styles: reuse existing style combos via copy/paste-special(formats)
#    so no redundant cellXfs entries get minted.
# =========================================================================

# C2:D7 -> same look as B (regular, black font, vertical-center)
$ws.Range("B2").Copy()
$ws.Range("C2:D7").PasteSpecial(-4122)

# D1:G1 -> same look as A1 (bold header)
$ws.Range("A1").Copy()
$ws.Range("D1:G1").PasteSpecial(-4122)

# =========================================================================
# 4) Column widths
# =========================================================================
$ws.Columns.Item(4).ColumnWidth = 22
$ws.Columns.Item(5).ColumnWidth = 39.1640625
$ws.Columns.Item(6).ColumnWidth = 39.1640625
$ws.Columns.Item(7).ColumnWidth = 42

# =========================================================================
# 5) Row heights
# =========================================================================
$ws.Rows.Item(3).RowHeight = 90
$ws.Rows.Item(4).RowHeight = 105
$ws.Rows.Item(5).RowHeight = 105
$ws.Rows.Item(6).RowHeight = 90
$ws.Rows.Item(7).RowHeight = 105

# =========================================================================
# 6) Sheet view: zoom + selection
# =========================================================================
$excel.ActiveWindow.Zoom = 150
$ws.Range("G2").Select()
